$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add default_position ("leggings") for the Thieves Leather row (row 77)
$ws.Range("G77").Value = "leggings"

# Row 243: Widow Basher (weapon)
$ws.Range("B243").Value = "Weapon Crafter Spell"
$ws.Range("C243").Value = 1
$ws.Range("D243").Value = "Widow Basher"
$ws.Range("E243").Value = "weapon"
$ws.Range("F243").Value = "A war-hammer nicknamed the Widow Basher. All because the guy who owned it, like to hunt down widows."
$ws.Range("H243").Value = 145
$ws.Range("K243").Value = 3600
$ws.Range("O243").Value = 0.2
$ws.Range("P243").Value = 0.2
$ws.Range("Q243").Value = 0.2
$ws.Range("R243").Value = 0.2
$ws.Range("S243").Value = 0.2
$ws.Range("U243").Value = 1
$ws.Range("V243").Value = 28
$ws.Range("W243").Value = 50
$ws.Range("X243").Value = "weapon"

# Row 244: Ripped Cloth (body)
$ws.Range("A244").Value = "Life Stealer"
$ws.Range("C244").Value = 1
$ws.Range("D244").Value = "Ripped Cloth"
$ws.Range("E244").Value = "body"
$ws.Range("F244").Value = "It's at least clothing, to say the least."
$ws.Range("G244").Value = "body"
$ws.Range("J244").Value = 4
$ws.Range("K244").Value = 10
$ws.Range("U244").Value = 1
$ws.Range("V244").Value = 1
$ws.Range("W244").Value = 5
$ws.Range("X244").Value = "armour"

# Row 245: Tin Helmet (helmet)
$ws.Range("B245").Value = "Weapons Glory"
$ws.Range("C245").Value = 1
$ws.Range("D245").Value = "Tin Helmet"
$ws.Range("E245").Value = "helmet"
$ws.Range("F245").Value = "Simple, not very sturdy, but simple."
$ws.Range("G245").Value = "helmet"
$ws.Range("J245").Value = 3
$ws.Range("K245").Value = 55
$ws.Range("U245").Value = 1
$ws.Range("V245").Value = 3
$ws.Range("W245").Value = 8
$ws.Range("X245").Value = "armour"

# Row 246: Thieves Leather (leggings)
$ws.Range("A246").Value = "Mages Inspiration"
$ws.Range("C246").Value = 1
$ws.Range("D246").Value = "Thieves Leather"
$ws.Range("E246").Value = "leggings"
$ws.Range("F246").Value = "Made from a leather that is blessed by a priest who was rumoured to be the greatest thief of all time."
$ws.Range("G246").Value = "leggings"
$ws.Range("J246").Value = 90
$ws.Range("K246").Value = 18000
$ws.Range("O246").Value = 0.14
$ws.Range("P246").Value = 0.14
$ws.Range("Q246").Value = 0.14
$ws.Range("R246").Value = 0.14
$ws.Range("S246").Value = 0.14
$ws.Range("U246").Value = 1
$ws.Range("V246").Value = 40
$ws.Range("W246").Value = 75
$ws.Range("X246").Value = "armour"

# Row 247: Witches Leggings (leggings)
$ws.Range("A247").Value = "Archbishops Prayer"
$ws.Range("C247").Value = 1
$ws.Range("D247").Value = "Witches Leggings"
$ws.Range("E247").Value = "leggings"
$ws.Range("F247").Value = "Not really sure what these are. A mixture of mesh and leather and bits of nature. I have no idea what these are but you wear them on your legs. Enjoy."
$ws.Range("G247").Value = "leggings"
$ws.Range("J247").Value = 5
$ws.Range("K247").Value = 100
$ws.Range("U247").Value = 1
$ws.Range("V247").Value = 6
$ws.Range("W247").Value = 12
$ws.Range("X247").Value = "armour"
